$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "AMSIN": append two new rows (31, 32) of registration-history data.
# ---------------------------------------------------------------------------
$wsAMSIN = $wb.Worksheets.Item("AMSIN")

# Insert two blank rows at the bottom of the table, inheriting the formatting
# of the row directly above (row 30), which is styled with style index 5
# (general) for A/C/D/E/F/G and the datetime style for B.
$wsAMSIN.Rows("31:32").Insert(-4121, 0)

# A scratch cell used to push literal text (not an auto-parsed date) into the
# date columns without disturbing the destination cell's inherited style.
$scratchIN = $wsAMSIN.Range("Z1")
$scratchIN.NumberFormat = "@"

$scratchIN.Value = "2023-02-07"
$scratchIN.Copy()
$wsAMSIN.Range("A31").PasteSpecial(-4163)

$scratchIN.Value = "2023-02-07"
$scratchIN.Copy()
$wsAMSIN.Range("A32").PasteSpecial(-4163)

$scratchIN.Clear()

$wsAMSIN.Range("B31").Value = 44964.49304700232
$wsAMSIN.Range("C31").Value = "testingpay172"
$wsAMSIN.Range("D31").Value = 41
$wsAMSIN.Range("E31").Value = 42
$wsAMSIN.Range("F31").Value = -1
$wsAMSIN.Range("G31").Value = 1.56

$wsAMSIN.Range("B32").Value = 44964.50708280093
$wsAMSIN.Range("C32").Value = "ppaytest1172"
$wsAMSIN.Range("D32").Value = 44
$wsAMSIN.Range("E32").Value = 43
$wsAMSIN.Range("F32").Value = 1
$wsAMSIN.Range("G32").Value = 2.3

# ---------------------------------------------------------------------------
# Sheet "AMS": restyle row 28, fix its run-time value, and append two new
# rows (29, 30) of registration-history data.
# ---------------------------------------------------------------------------
$wsAMS = $wb.Worksheets.Item("AMS")

# Insert two blank rows at position 28, inheriting the style-5/style-10
# formatting of row 27 above. This pushes the existing (unstyled) row 28
# ("pay172three") down to row 30, which conveniently matches the target
# layout: the final row 30 should remain unstyled.
$wsAMS.Rows("28:29").Insert(-4121, 0)

$scratchAMS = $wsAMS.Range("Z1")
$scratchAMS.NumberFormat = "@"

$scratchAMS.Value = "2023-02-03"
$scratchAMS.Copy()
$wsAMS.Range("A28").PasteSpecial(-4163)

$scratchAMS.Value = "2023-02-07"
$scratchAMS.Copy()
$wsAMS.Range("A29").PasteSpecial(-4163)

$scratchAMS.Value = "2023-02-07"
$scratchAMS.Copy()
$wsAMS.Range("A30").PasteSpecial(-4163)

$scratchAMS.Clear()

# Row 28 ("pay172three") now newly styled, with the corrected run-time value.
$wsAMS.Range("B28").Value = 44960.63808569445
$wsAMS.Range("C28").Value = "pay172three"
$wsAMS.Range("D28").Value = 41
$wsAMS.Range("E28").Value = 39
$wsAMS.Range("F28").Value = 2
$wsAMS.Range("G28").Value = 1.03

# Row 29 ("testingpay172"), newly added, styled like row 28.
$wsAMS.Range("B29").Value = 44964.5041196875
$wsAMS.Range("C29").Value = "testingpay172"
$wsAMS.Range("D29").Value = 44
$wsAMS.Range("E29").Value = 43
$wsAMS.Range("F29").Value = 1
$wsAMS.Range("G29").Value = 1.56

# Row 30 ("testff172"), newly added, left unstyled (matches the original
# row 28 formatting that shifted down to this position).
$wsAMS.Range("B30").Value = 44964.57731602514
$wsAMS.Range("C30").Value = "testff172"
$wsAMS.Range("D30").Value = 44
$wsAMS.Range("E30").Value = 44
$wsAMS.Range("F30").Value = 0
$wsAMS.Range("G30").Value = 1.36
